$d = $word.ActiveDocument
$vt = [char]11   # vertical-tab == <w:br/> line break inside a run

# ---------------------------------------------------------------------
# 1) Requisitos list: drop the "LOB1009 - Leitura e Interpretação de
#    Desenho Técnico (Requisito)" line entirely (its text run plus the
#    line-break that terminates it).
# ---------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute( `
    "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $doomed = $d.Range($find1.Start, $find1.End + 1)
    $doomed.Delete()
}

# ---------------------------------------------------------------------
# 2) Requisitos list: the three "LOQ..." chemistry lines are replaced by
#    four new lines (Química Geral I/II retired, Experimental renamed,
#    two "Fundamentos" courses and a CAD course added). Rather than
#    rewriting the existing runs' text in place (which coalesces
#    neighbouring same-format runs into one <w:r>), delete the whole
#    old block and insert four fresh runs so each line keeps its own
#    <w:r><w:t/><w:br/></w:r>.
# ---------------------------------------------------------------------
$findStart = $d.Content
$null = $findStart.Find.Execute( `
    "LOQ4031 -  Química Geral I  (Requisito)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockStart = $findStart.Start

$findEnd = $d.Content
$null = $findEnd.Find.Execute( `
    "LOQ4095 -  Química Geral Experimental  (Requisito)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockEnd = $findEnd.End + 1

$oldBlock = $d.Range($blockStart, $blockEnd)
$oldBlock.Delete()

$insertionPoint = $d.Range($blockStart, $blockStart)
$insertionPoint.InsertAfter("LOQ4095 -  Química Geral Experimental  (Requisito)" + $vt)
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter("LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)" + $vt)
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter("LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)" + $vt)
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter("LOQ4247 -  Desenho Assistido por Computador  (Requisito)" + $vt)
$insertionPoint.Collapse(0)
